# Applies the cryptos price/volume refresh described by the commit:
# "Updated cryptos list on Tue Oct 10 09:43:18 UTC 2023 with GitHub Actions"
#
# For each changed cell we briefly force a Text number format so that
# numeric-looking strings (e.g. "209.40", "0.0868") are kept as literal
# text instead of being parsed into doubles (which would drop trailing
# zeros / switch to scientific notation). We then restore the cell style
# to "Normal" so no stray formatting is left behind, matching the
# original (unstyled) D/E data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.747.44'
Set-TextValue $ws.Range("E2") '  +0.01%  '
Set-TextValue $ws.Range("D3") '1.594.70'
Set-TextValue $ws.Range("E3") '  -1.27%  '
Set-TextValue $ws.Range("E4") '  +0.19%  '
Set-TextValue $ws.Range("D5") '209.40'
Set-TextValue $ws.Range("E5") '  -0.61%  '
Set-TextValue $ws.Range("E6") '  -1.99%  '
Set-TextValue $ws.Range("E7") '  +0.27%  '
Set-TextValue $ws.Range("D8") '22.35'
Set-TextValue $ws.Range("E8") '  -2.23%  '
Set-TextValue $ws.Range("E9") '  -1.30%  '
Set-TextValue $ws.Range("E10") '  -1.50%  '
Set-TextValue $ws.Range("D11") '0.0868'
Set-TextValue $ws.Range("E11") '  -1.28%  '
Set-TextValue $ws.Range("D12") '1.821.56'
Set-TextValue $ws.Range("E12") '  -1.32%  '
Set-TextValue $ws.Range("D13") '1.598.54'
Set-TextValue $ws.Range("E13") '  -0.95%  '
Set-TextValue $ws.Range("E14") '  -2.20%  '
Set-TextValue $ws.Range("D15") '0.532'
Set-TextValue $ws.Range("E15") '  -3.33%  '
Set-TextValue $ws.Range("D16") '27.757.89'
Set-TextValue $ws.Range("E16") '  +0.01%  '
Set-TextValue $ws.Range("D17") '63.43'
Set-TextValue $ws.Range("E17") '  -1.43%  '
Set-TextValue $ws.Range("D18") '219.30'
Set-TextValue $ws.Range("E18") '  -2.75%  '
Set-TextValue $ws.Range("D19") '7.37'
Set-TextValue $ws.Range("E19") '  -2.40%  '
Set-TextValue $ws.Range("E20") '  -1.84%  '
Set-TextValue $ws.Range("E22") '  -2.97%  '
Set-TextValue $ws.Range("D23") '9.74'
Set-TextValue $ws.Range("E23") '  -1.95%  '
Set-TextValue $ws.Range("E24") '  -3.38%  '
Set-TextValue $ws.Range("D25") '153.88'
Set-TextValue $ws.Range("E25") '  -0.51%  '
Set-TextValue $ws.Range("D26") '7.13'
Set-TextValue $ws.Range("E26") '  +3.90%  '
Set-TextValue $ws.Range("E27") '  +0.24%  '
Set-TextValue $ws.Range("D28") '15.16'
Set-TextValue $ws.Range("E28") '  -0.68%  '
Set-TextValue $ws.Range("E29") '  -2.85%  '
Set-TextValue $ws.Range("E30") '  -0.76%  '
Set-TextValue $ws.Range("E31") '  -0.77%  '
Set-TextValue $ws.Range("E32") '  -3.68%  '
Set-TextValue $ws.Range("D33") '1.381.26'
Set-TextValue $ws.Range("E33") '  -0.75%  '
Set-TextValue $ws.Range("D34") '2.98'
Set-TextValue $ws.Range("E34") '  -2.15%  '
Set-TextValue $ws.Range("E35") '  -3.02%  '
Set-TextValue $ws.Range("D36") '0.974'
Set-TextValue $ws.Range("E36") '  +0.23%  '
Set-TextValue $ws.Range("E38") '  +0.58%  '
Set-TextValue $ws.Range("D39") '0.537'
Set-TextValue $ws.Range("E39") '  -2.35%  '
Set-TextValue $ws.Range("D40") '0.828'
Set-TextValue $ws.Range("E40") '  -1.41%  '
Set-TextValue $ws.Range("E41") '  +0.22%  '
Set-TextValue $ws.Range("E42") '  -1.91%  '
Set-TextValue $ws.Range("D43") '64.53'
Set-TextValue $ws.Range("E43") '  -0.69%  '
Set-TextValue $ws.Range("E44") '  +2.71%  '
Set-TextValue $ws.Range("E45") '  -0.44%  '
Set-TextValue $ws.Range("E46") '  -2.10%  '
Set-TextValue $ws.Range("D47") '1.732.13'
Set-TextValue $ws.Range("E47") '  -1.39%  '
Set-TextValue $ws.Range("D48") '86.15'
Set-TextValue $ws.Range("E48") '  -3.55%  '
Set-TextValue $ws.Range("E49") '  +0.46%  '
Set-TextValue $ws.Range("D50") '0.0965'
Set-TextValue $ws.Range("E50") '  -2.27%  '
Set-TextValue $ws.Range("D51") '0.0495'
Set-TextValue $ws.Range("E51") '  -1.13%  '
